$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: remove all existing content/formatting so the sheet can be
# rebuilt to match the new "sample protocol & manual" layout.
$ws.Cells.Clear()

# --- Header / comment lines (rows 1-8) ---
$ws.Range("A1").Value = "* sample protocol file"
$ws.Range("A2").Value = "* "
$ws.Range("A3").Value = "*"
$ws.Range("B3").Value = 'line that begins with "*" is regarded as comment, and will be ignored'
$ws.Range("A4").Value = "*"

$ws.Range("A5").Value = "*"
$ws.Range("B5").Value = 'The Column A must be "Trial"'
$ws.Range("A6").Value = "*"
$ws.Range("B6").Value = "Experimental parameters should be set in Column B, C, …"
$ws.Range("A7").Value = "*"
$ws.Range("B7").Value = "The fied names must match those used in config.csv"

$ws.Range("A8").Value = "*"

# Highlight the manual/instructions block (rows 5-7) in yellow.
$ws.Range("A5:B7").Interior.Color = 65535

# --- Table header (row 9) ---
$ws.Range("A9").Value = "Trial"
$ws.Range("B9").Value = "Target Direction"
$ws.Range("C9").Value = "Visuo Motor Rotation Angle"
$ws.Range("D9").Value = "Is Mirror Reversed"
$ws.Range("E9").Value = "Vibration"
$ws.Range("F9").Value = "Start Area Position X"

# Shade the table header row using theme "Background 2".
$ws.Range("A9:F9").Interior.Color = 65535
$ws.Range("A9:F9").Interior.ThemeColor = 4

# --- Table data (rows 10-14) ---
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 80
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = $true
$ws.Range("F10").Value = 0

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 100
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = $false
$ws.Range("E11").Value = $true
$ws.Range("F11").Value = 0

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = $false
$ws.Range("E12").Value = $true
$ws.Range("F12").Value = 0

$ws.Range("A13").Value = 4
$ws.Range("B13").Value = 80
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = 0

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = 100
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = $true
$ws.Range("E14").Value = $false
$ws.Range("F14").Value = 0

# Stray leftover cell in row 15.
$ws.Range("F15").Value = 0

# --- Column widths (B:F) ---
$ws.Columns("B:F").ColumnWidth = 22.5

# --- Selection cursor ---
$ws.Range("E18").Select() | Out-Null
